$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'301.27"
$ws.Range("E2").Value = "'0.72%"
$ws.Range("D3").Value = "'31.64"
$ws.Range("E3").Value = "'1.36%"
$ws.Range("D4").Value = "'5.097"
$ws.Range("E4").Value = "'-1.28%"
$ws.Range("E5").Value = "'-3.02%"
$ws.Range("D6").Value = "'2.237"
$ws.Range("E6").Value = "'-16.61%"
$ws.Range("D7").Value = "'7.796"
$ws.Range("E7").Value = "'-0.33%"
$ws.Range("D8").Value = "'3.827"
$ws.Range("E8").Value = "'0.05%"
$ws.Range("D9").Value = "'0.9167"
$ws.Range("E9").Value = "'0.26%"
$ws.Range("D10").Value = "'0.1756"
$ws.Range("E10").Value = "'1.01%"
$ws.Range("D11").Value = "'0.07550"
$ws.Range("E11").Value = "'4.03%"
$ws.Range("D12").Value = "'0.09021"
$ws.Range("E12").Value = "'7.86%"
$ws.Range("D13").Value = "'0.03088"
$ws.Range("E13").Value = "'3.19%"
$ws.Range("E14").Value = "'0.68%"
$ws.Range("D15").Value = "'0.001512"
$ws.Range("E15").Value = "'1.48%"
$ws.Range("D16").Value = "'0.006044"
$ws.Range("E16").Value = "'0.17%"
$ws.Range("E17").Value = "'-0.85%"
$ws.Range("D18").Value = "'2.252"
$ws.Range("E18").Value = "'0.04%"
$ws.Range("D19").Value = "'0.3268"
$ws.Range("E19").Value = "'-0.49%"
$ws.Range("E20").Value = "'0.68%"
$ws.Range("D21").Value = "'4.334"
$ws.Range("E21").Value = "'-6.60%"
$ws.Range("D23").Value = "'0.04584"
$ws.Range("E23").Value = "'0.25%"
$ws.Range("D24").Value = "'0.001252"
$ws.Range("E24").Value = "'-0.48%"
$ws.Range("D25").Value = "'0.004463"
$ws.Range("E25").Value = "'0.24%"
$ws.Range("E26").Value = "'5.88%"
$ws.Range("E27").Value = "'-1.40%"
$ws.Range("D39").Value = "'0.01768"
$ws.Range("E39").Value = "'-4.04%"
$ws.Range("D40").Value = "'0.04810"
$ws.Range("E40").Value = "'6.38%"
$ws.Range("D41").Value = "'0.007511"
$ws.Range("E41").Value = "'6.52%"
$ws.Range("D42").Value = "'0.1357"
$ws.Range("E42").Value = "'1.02%"
$ws.Range("D43").Value = "'0.002188"
$ws.Range("E43").Value = "'-2.28%"
$ws.Range("D44").Value = "'0.01023"
$ws.Range("E44").Value = "'-2.80%"
$ws.Range("D45").Value = "'0.00006191"
$ws.Range("E45").Value = "'-4.56%"
$ws.Range("E46").Value = "'-0.01%"
$ws.Range("E47").Value = "'28.81%"
$ws.Range("D48").Value = "'1.146"
$ws.Range("E48").Value = "'39.68%"
$ws.Range("D49").Value = "'0.00002100"
$ws.Range("E49").Value = "'-0.01%"
$ws.Range("D50").Value = "'0.0002000"
$ws.Range("E50").Value = "'-0.01%"
